# fix: mejora en la exportacion de detalles form documento
#
# The "TIPOS-LISTAS" helper sheet used to feed the "Tipo de Documento"
# dropdown had two entries removed ("07 - Nota de credito" and
# "08 - Nota de debito"), so every list below them shifts up two rows
# and the used range shrinks from D2:D20 to D2:D18.

$wb = $excel.ActiveWorkbook

# Remember which sheet/cell was active so we can restore it after we
# touch the helper sheet's view (zoom/selection) further down.
$originalActive = $wb.ActiveSheet

$wsListas = $wb.Worksheets.Item("TIPOS-LISTAS")

# Remove the two obsolete document-type rows (D8:D9) and shift the
# remaining entries up, which also drops the now-unused shared strings
# and shrinks the sheet's dimension from B1:H20 to B1:H18.
$wsListas.Range("D8:D9").Delete(-4162)  # xlShiftUp

# The defined name that backs the "Tipo de Documento" validation list
# needs to track the new, shorter range.
$wb.Names.Item("TIPODOCUMENTO").RefersTo = "='TIPOS-LISTAS'!`$D`$2:`$D`$18"

# Column D on TIPOS-LISTAS was widened slightly by the author while
# reviewing the shortened list.
$wsListas.Columns.Item(4).ColumnWidth = 40.5

# Reflect the zoom level / active cell the author left the helper sheet
# on, then switch back to whatever sheet was active before.
$wsListas.Select()
$wsListas.Range("D5").Select()
$excel.ActiveWindow.Zoom = 85
$originalActive.Select()
